# Scheduled runner update: refresh market-price-derived profit figures
# (currentAveragePrice*, Leve*Price*, LeveProfit*) across the per-job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 234.91667
$ws.Range("I9").Value = 155
$ws.Range("J9").Value = 302.53845
$ws.Range("K9").Value = 155
$ws.Range("L9").Value = 302.53845
$ws.Range("M9").Value = 14
$ws.Range("N9").Value = -640.53845
$ws.Range("H33").Value = 224.23077
$ws.Range("I33").Value = 248.1
$ws.Range("J33").Value = 144.66667
$ws.Range("K33").Value = 248.1
$ws.Range("L33").Value = 144.66667
$ws.Range("M33").Value = -19.09999999999999
$ws.Range("N33").Value = -602.6666700000001
$ws.Range("H69").Value = 51979.4
$ws.Range("I69").Value = 1500
$ws.Range("J69").Value = 64599.25
$ws.Range("K69").Value = 4500
$ws.Range("L69").Value = 193797.75
$ws.Range("M69").Value = -3626
$ws.Range("N69").Value = -195545.75
$ws.Range("H72").Value = 51979.4
$ws.Range("I72").Value = 1500
$ws.Range("J72").Value = 64599.25
$ws.Range("K72").Value = 13500
$ws.Range("L72").Value = 581393.25
$ws.Range("M72").Value = -9132
$ws.Range("N72").Value = -590129.25
$ws.Range("H74").Value = 9432.5
$ws.Range("I74").Value = 8790.625
$ws.Range("J74").Value = 12000
$ws.Range("K74").Value = 8790.625
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = -7854.625
$ws.Range("N74").Value = -13872
$ws.Range("H77").Value = 9432.5
$ws.Range("I77").Value = 8790.625
$ws.Range("J77").Value = 12000
$ws.Range("K77").Value = 43953.125
$ws.Range("L77").Value = 60000
$ws.Range("M77").Value = -39273.125
$ws.Range("N77").Value = -69360
$ws.Range("H80").Value = 2970.138
$ws.Range("I80").Value = 1248.1666
$ws.Range("J80").Value = 4185.647
$ws.Range("K80").Value = 3744.4998
$ws.Range("L80").Value = 12556.941
$ws.Range("M80").Value = -2746.4998
$ws.Range("N80").Value = -14552.941
$ws.Range("H82").Value = 3749
$ws.Range("I82").Value = 3749
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 11247
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 2970.138
$ws.Range("I83").Value = 1248.1666
$ws.Range("J83").Value = 4185.647
$ws.Range("K83").Value = 11233.4994
$ws.Range("L83").Value = 37670.823
$ws.Range("M83").Value = -6241.499400000001
$ws.Range("N83").Value = -47654.823
$ws.Range("H85").Value = 3749
$ws.Range("I85").Value = 3749
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 11247
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 2519.647
$ws.Range("J86").Value = 3199.625
$ws.Range("L86").Value = 3199.625
$ws.Range("N86").Value = -5445.625
$ws.Range("H89").Value = 2519.647
$ws.Range("J89").Value = 3199.625
$ws.Range("L89").Value = 15998.125
$ws.Range("N89").Value = -27230.125
$ws.Range("H132").Value = 21333.822
$ws.Range("I132").Value = 2083.2856
$ws.Range("K132").Value = 6249.8568
$ws.Range("M132").Value = -3719.8568
$ws.Range("H137").Value = 4657131
$ws.Range("I137").Value = 6674001
$ws.Range("J137").Value = 2815
$ws.Range("K137").Value = 20022003
$ws.Range("L137").Value = 8445
$ws.Range("M137").Value = -20019453
$ws.Range("N137").Value = -13545
$ws.Range("H138").Value = 7328.186
$ws.Range("I138").Value = 3138.182
$ws.Range("J138").Value = 8768.5
$ws.Range("K138").Value = 9414.545999999998
$ws.Range("L138").Value = 26305.5
$ws.Range("M138").Value = -4274.545999999998
$ws.Range("N138").Value = -36585.5

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14090578
$ws.Range("I32").Value = 15390989
$ws.Range("J32").Value = 2787.1667
$ws.Range("K32").Value = 15390989
$ws.Range("L32").Value = 2787.1667
$ws.Range("M32").Value = -15390702
$ws.Range("N32").Value = -3361.1667
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50976
$ws.Range("H45").Value = 4796.2
$ws.Range("I45").Value = 3184
$ws.Range("K45").Value = 3184
$ws.Range("M45").Value = -2807
$ws.Range("H61").Value = 4003.926
$ws.Range("I61").Value = 2631.1667
$ws.Range("K61").Value = 2631.1667
$ws.Range("M61").Value = -2419.1667
$ws.Range("H119").Value = 79698
$ws.Range("J119").Value = 79698
$ws.Range("L119").Value = 79698
$ws.Range("N119").Value = -89374
$ws.Range("H136").Value = 4003.926
$ws.Range("I136").Value = 2631.1667
$ws.Range("K136").Value = 7893.500100000001
$ws.Range("M136").Value = -5343.500100000001

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5641.9473
$ws.Range("I107").Value = 5293.5625
$ws.Range("K107").Value = 5293.5625
$ws.Range("M107").Value = -3373.5625
$ws.Range("H134").Value = 3000.8071
$ws.Range("I134").Value = 2822.0625
$ws.Range("J134").Value = 3229.6
$ws.Range("K134").Value = 8466.1875
$ws.Range("L134").Value = 9688.799999999999
$ws.Range("M134").Value = -5931.1875
$ws.Range("N134").Value = -14758.8
$ws.Range("H135").Value = 39800
$ws.Range("J135").Value = 39800
$ws.Range("L135").Value = 39800
$ws.Range("N135").Value = -49940

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2644.1406
$ws.Range("I31").Value = 1808.4783
$ws.Range("J31").Value = 4779.722
$ws.Range("K31").Value = 1808.4783
$ws.Range("L31").Value = 4779.722
$ws.Range("M31").Value = -1513.4783
$ws.Range("N31").Value = -5369.722
$ws.Range("H34").Value = 2644.1406
$ws.Range("I34").Value = 1808.4783
$ws.Range("J34").Value = 4779.722
$ws.Range("K34").Value = 1808.4783
$ws.Range("L34").Value = 4779.722
$ws.Range("M34").Value = -1606.4783
$ws.Range("N34").Value = -5183.722
$ws.Range("H58").Value = 2029.5385
$ws.Range("I58").Value = 1661.5714
$ws.Range("J58").Value = 3575
$ws.Range("K58").Value = 1661.5714
$ws.Range("L58").Value = 3575
$ws.Range("M58").Value = -1458.5714
$ws.Range("N58").Value = -3981
$ws.Range("H118").Value = 34250
$ws.Range("J118").Value = 34250
$ws.Range("L118").Value = 34250
$ws.Range("N118").Value = -37564
$ws.Range("H132").Value = 3605.9583
$ws.Range("I132").Value = 3605.9583
$ws.Range("K132").Value = 10817.8749
$ws.Range("M132").Value = -8287.874899999999
$ws.Range("H134").Value = 3311.8333
$ws.Range("I134").Value = 3311.8333
$ws.Range("K134").Value = 9935.499899999999
$ws.Range("M134").Value = -7400.499899999999
$ws.Range("H136").Value = 2029.5385
$ws.Range("I136").Value = 1661.5714
$ws.Range("J136").Value = 3575
$ws.Range("K136").Value = 4984.7142
$ws.Range("L136").Value = 10725
$ws.Range("M136").Value = -2434.7142
$ws.Range("N136").Value = -15825

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 649.5833
$ws.Range("I92").Value = 592.1429000000001
$ws.Range("J92").Value = 730
$ws.Range("K92").Value = 1776.4287
$ws.Range("L92").Value = 2190
$ws.Range("M92").Value = -528.4287000000002
$ws.Range("N92").Value = -4686

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 27880
$ws.Range("I46").Value = 23000
$ws.Range("J46").Value = 31133.334
$ws.Range("K46").Value = 23000
$ws.Range("L46").Value = 31133.334
$ws.Range("M46").Value = -22844
$ws.Range("N46").Value = -31445.334
$ws.Range("H57").Value = 15000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 15000
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -16640
$ws.Range("H132").Value = 6330.07
$ws.Range("I132").Value = 5690.3516
$ws.Range("J132").Value = 10275
$ws.Range("K132").Value = 17071.0548
$ws.Range("L132").Value = 30825
$ws.Range("M132").Value = -14541.0548
$ws.Range("N132").Value = -35885

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2219.394
$ws.Range("I46").Value = 974.9
$ws.Range("J46").Value = 2760.4783
$ws.Range("K46").Value = 974.9
$ws.Range("L46").Value = 2760.4783
$ws.Range("M46").Value = -786.9
$ws.Range("N46").Value = -3136.4783
$ws.Range("H132").Value = 5431.304
$ws.Range("I132").Value = 5088.9697
$ws.Range("J132").Value = 6300.3076
$ws.Range("K132").Value = 15266.9091
$ws.Range("L132").Value = 18900.9228
$ws.Range("M132").Value = -12736.9091
$ws.Range("N132").Value = -23960.9228
$ws.Range("H136").Value = 3709
$ws.Range("I136").Value = 3216.1667
$ws.Range("J136").Value = 6666
$ws.Range("K136").Value = 9648.500100000001
$ws.Range("L136").Value = 19998
$ws.Range("M136").Value = -7098.500100000001
$ws.Range("N136").Value = -25098

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 77523.336
$ws.Range("J119").Value = 77523.336
$ws.Range("L119").Value = 77523.336
$ws.Range("N119").Value = -87199.336
$ws.Range("H136").Value = 9121.777
$ws.Range("I136").Value = 9506.571
$ws.Range("J136").Value = 7775
$ws.Range("K136").Value = 28519.713
$ws.Range("L136").Value = 23325
$ws.Range("M136").Value = -25969.713
$ws.Range("N136").Value = -28425
$ws.Range("H141").Value = 84268.27
$ws.Range("J141").Value = 86630.10000000001
$ws.Range("L141").Value = 86630.10000000001
$ws.Range("N141").Value = -96990.10000000001

